# Task.xlsx edit: "unify the conception of DataNode, DataTable, Entity."
#
# The sheet that used to represent a single "Property" table is renamed to
# "DataNode" so the naming matches the unified DataNode/DataTable/Entity
# concept used across the rest of the config. Along with the rename, the
# editor had cleared the manual row-height override on the header/sub-header
# row (row 7) and left the cursor resting on F25 in the frozen (data) pane
# before saving.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Rename the worksheet: "Property1" -> "DataNode"
$ws.Name = "DataNode"

# 2. Clear the explicit row height on row 7 (the "TRUE/FALSE" validation
#    row), restoring it to the sheet's default row height instead of the
#    old manually-set 14pt.
$ws.Rows.Item(7).AutoFit()

# 3. Leave the selection on F25 within the frozen lower-left pane, matching
#    the cursor position the file was saved with.
$ws.Range("F25").Select()
